$wb = $excel.ActiveWorkbook

# ======================================================================
# 1) "总计" (summary) sheet: insert a new "2022-Q4" row right after the
#    header, pushing the existing quarters down by one row.
# ======================================================================
$summary = $wb.Worksheets.Item(1)

# Shift rows 2-6 down to 3-7 (Excel-native row insert keeps formatting
# of the shifted rows intact, including the bordered/bold index column).
$summary.Rows(2).Insert()

# The newly inserted blank row inherited the header's border on B:D -
# the target row has no special formatting there, so clear it.
$summary.Range("B2:D2").ClearFormats()

# Give the new index cell (A2) the same bold/border/center style used by
# the rest of the index column (copy format only, not value).
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

# Fill the new 2022-Q4 row.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.09

# Keep the index column (A) a simple 0..5 running sequence.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5

# ======================================================================
# 2) Brand-new "2022-Q4" worksheet, inserted right before "2022-Q2" so
#    sheet order becomes 总计, 2022-Q4, 2022-Q2, 2022-Q1, 2021-Q3,
#    2021-Q2, 2021-Q1.
# ======================================================================
$insertBefore = $wb.Worksheets.Item(2)
$q4Sheet = $wb.Worksheets.Add($insertBefore)
$q4Sheet.Name = "2022-Q4"

# NOTE: fetch "2022-Q2" by name (not by the old positional index) - once
# the new sheet is inserted before it, any Item(2)-style reference grabbed
# beforehand now resolves to the newly inserted sheet instead.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Match the header-row formatting used on every other quarter sheet.
$q2Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

# Match the index-column (A) formatting on data rows too.
$q2Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A4").PasteSpecial(-4122)

# Header row.
$q4Sheet.Cells.Item(1,2).Value = "基金代码"
$q4Sheet.Cells.Item(1,3).Value = "基金名称"
$q4Sheet.Cells.Item(1,4).Value = "基金规模"
$q4Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1,6).Value = "仓位占比"
$q4Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1,8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percents)
# that must stay text, same as the other quarter sheets - force text
# format before writing so Excel doesn't silently coerce them to numbers.
# (C is the fund name, always non-numeric, so it's left alone.)
$q4Sheet.Range("B2:B4").NumberFormat = "@"
$q4Sheet.Range("D2:G4").NumberFormat = "@"

# Row 2 - 970007
$q4Sheet.Cells.Item(2,1).Value = 0
$q4Sheet.Cells.Item(2,2).Value = "970007"
$q4Sheet.Cells.Item(2,3).Value = "华安证券汇赢增利一年持有混合B"
$q4Sheet.Cells.Item(2,4).Value = "11.05"
$q4Sheet.Cells.Item(2,5).Value = "22.39"
$q4Sheet.Cells.Item(2,6).Value = "0.48"
$q4Sheet.Cells.Item(2,7).Value = "0.0530"
$q4Sheet.Cells.Item(2,8).Value = 9

# Row 3 - 970008
$q4Sheet.Cells.Item(3,1).Value = 1
$q4Sheet.Cells.Item(3,2).Value = "970008"
$q4Sheet.Cells.Item(3,3).Value = "华安证券汇赢增利一年持有混合C"
$q4Sheet.Cells.Item(3,4).Value = "8.56"
$q4Sheet.Cells.Item(3,5).Value = "22.39"
$q4Sheet.Cells.Item(3,6).Value = "0.48"
$q4Sheet.Cells.Item(3,7).Value = "0.0411"
$q4Sheet.Cells.Item(3,8).Value = 9

# Row 4 - 970006
$q4Sheet.Cells.Item(4,1).Value = 2
$q4Sheet.Cells.Item(4,2).Value = "970006"
$q4Sheet.Cells.Item(4,3).Value = "华安证券汇赢增利一年持有混合A"
$q4Sheet.Cells.Item(4,4).Value = "0.18"
$q4Sheet.Cells.Item(4,5).Value = "22.39"
$q4Sheet.Cells.Item(4,6).Value = "0.48"
$q4Sheet.Cells.Item(4,7).Value = "0.0009"
$q4Sheet.Cells.Item(4,8).Value = 9
